$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "ATA_WEAPON_CLASS_HORRIFIC"
$ws.Range("B15").Value = "Horrific"
$ws.Range("C15").Value = "暗渊"

$ws.Range("C15").Select() | Out-Null
